$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = $false
$ws.Range("AG2").Value = "Fields missing .Couldn't Verify Your aadhar card."

$ws.Range("AF3").Value = $true

$ws.Range("AF4").Value = $true

$ws.Range("AF5").Value = $false
